$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $val = $cell.Value()
    if ($val -ne $null -and $val -is [string] -and $val.Contains("Pfizer Monovalent")) {
        $cell.Value = $val.Replace("Pfizer Monovalent", "Pfizer_mono")
    }
}
